$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 215
$ws.Range("I2").Value = 618
$ws.Range("J2").Value = 2509
$ws.Range("L2").Value = 708
$ws.Range("M2").Value = 43
$ws.Range("N2").Value = 436
$ws.Range("O2").Value = 3
$ws.Range("Q2").Value = 4
$ws.Range("R2").Value = 36
$ws.Range("S2").Value = 297
$ws.Range("T2").Value = 452
$ws.Range("V2").Value = 4004
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 3905
$ws.Range("Y2").Value = 9
$ws.Range("Z2").Value = 50
$ws.Range("AA2").Value = 23
